# Update countries & provincias Spain
# Applies the 6-Apr-2020 07:22 data refresh: updated totals for several
# countries plus the resulting re-sort (row text swaps) and the refreshed
# "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 6 de Abril de 2020 a las 07:22'

# Row 42
$ws.Cells.Item(42, 2).Value = 2220
$ws.Cells.Item(42, 3).Value = 51
$ws.Cells.Item(42, 5).Value = 1401
$ws.Cells.Item(42, 7).Value = 3
$ws.Cells.Item(42, 8).Value = 26

# Row 43
$ws.Cells.Item(43, 6).Value = 293

# Row 68
$ws.Cells.Item(68, 5).Value = 790
$ws.Cells.Item(68, 7).Value = 1
$ws.Cells.Item(68, 8).Value = 14

# Row 69
$ws.Cells.Item(69, 2).Value = 744
$ws.Cells.Item(69, 3).Value = 11
$ws.Cells.Item(69, 4).Value = 67
$ws.Cells.Item(69, 5).Value = 639
$ws.Cells.Item(69, 7).Value = 4
$ws.Cells.Item(69, 8).Value = 38

# Row 80
$ws.Cells.Item(80, 1).Value = 'Bulgaria'
$ws.Cells.Item(80, 2).Value = 541
$ws.Cells.Item(80, 3).Value = 10
$ws.Cells.Item(80, 4).Value = 37
$ws.Cells.Item(80, 5).Value = 483
$ws.Cells.Item(80, 6).Value = 22
$ws.Cells.Item(80, 7).Value = 1
$ws.Cells.Item(80, 8).Value = 21

# Row 81
$ws.Cells.Item(81, 1).Value = 'Letonia'
$ws.Cells.Item(81, 2).Value = 533
$ws.Cells.Item(81, 4).Value = 1
$ws.Cells.Item(81, 5).Value = 531
$ws.Cells.Item(81, 6).Value = 4
$ws.Cells.Item(81, 8).Value = 1

# Row 88
$ws.Cells.Item(88, 1).Value = 'Uzbekistan'
$ws.Cells.Item(88, 2).Value = 390
$ws.Cells.Item(88, 3).Value = 48
$ws.Cells.Item(88, 4).Value = 30
$ws.Cells.Item(88, 5).Value = 358
$ws.Cells.Item(88, 6).Value = 8
$ws.Cells.Item(88, 8).Value = 2

# Row 89
$ws.Cells.Item(89, 1).Value = 'Afganistan'
$ws.Cells.Item(89, 2).Value = 367
$ws.Cells.Item(89, 3).Value = 18
$ws.Cells.Item(89, 4).Value = 17
$ws.Cells.Item(89, 5).Value = 343
$ws.Cells.Item(89, 8).Value = 7

# Row 90
$ws.Cells.Item(90, 1).Value = 'Taiwan'
$ws.Cells.Item(90, 2).Value = 363
$ws.Cells.Item(90, 4).Value = 54
$ws.Cells.Item(90, 5).Value = 304
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 8).Value = 5

# Row 91
$ws.Cells.Item(91, 1).Value = 'Albania'
$ws.Cells.Item(91, 2).Value = 361
$ws.Cells.Item(91, 4).Value = 104
$ws.Cells.Item(91, 5).Value = 237
$ws.Cells.Item(91, 6).Value = 7
$ws.Cells.Item(91, 8).Value = 20

# Row 92
$ws.Cells.Item(92, 1).Value = 'Burkina Faso'
$ws.Cells.Item(92, 4).Value = 90
$ws.Cells.Item(92, 5).Value = 238
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 8).Value = 17

# Row 93
$ws.Cells.Item(93, 1).Value = 'Jordania'
$ws.Cells.Item(93, 2).Value = 345
$ws.Cells.Item(93, 4).Value = 110
$ws.Cells.Item(93, 5).Value = 230
$ws.Cells.Item(93, 6).Value = 5
$ws.Cells.Item(93, 8).Value = 5

# Row 94
$ws.Cells.Item(94, 1).Value = 'Reunion'
$ws.Cells.Item(94, 2).Value = 344
$ws.Cells.Item(94, 4).Value = 40
$ws.Cells.Item(94, 5).Value = 304
$ws.Cells.Item(94, 6).Value = 4
$ws.Cells.Item(94, 8).Value = 0

# Row 100
$ws.Cells.Item(100, 4).Value = 91
$ws.Cells.Item(100, 5).Value = 150

# Row 123
$ws.Cells.Item(123, 4).Value = 53
$ws.Cells.Item(123, 5).Value = 61

# Row 124
$ws.Cells.Item(124, 6).Value = 8

# Row 149
$ws.Cells.Item(149, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(149, 3).Value = 12
$ws.Cells.Item(149, 4).Value = 1
$ws.Cells.Item(149, 5).Value = 30
$ws.Cells.Item(149, 7).Value = 2
$ws.Cells.Item(149, 8).Value = 6

# Row 150
$ws.Cells.Item(150, 1).Value = 'Bermudas'
$ws.Cells.Item(150, 2).Value = 37
$ws.Cells.Item(150, 4).Value = 14
$ws.Cells.Item(150, 5).Value = 23
$ws.Cells.Item(150, 8).Value = 0

# Row 151
$ws.Cells.Item(151, 1).Value = 'Guam'
$ws.Cells.Item(151, 4).Value = 0
$ws.Cells.Item(151, 5).Value = 31
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 8).Value = 1

# Row 152
$ws.Cells.Item(152, 1).Value = 'San Martin (Parte Francesa)'
$ws.Cells.Item(152, 2).Value = 32
$ws.Cells.Item(152, 4).Value = 7
$ws.Cells.Item(152, 5).Value = 23
$ws.Cells.Item(152, 6).Value = 6
$ws.Cells.Item(152, 8).Value = 2

# Row 153
$ws.Cells.Item(153, 1).Value = 'Eritrea'
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(153, 5).Value = 29
$ws.Cells.Item(153, 8).Value = 0

# Row 154
$ws.Cells.Item(154, 1).Value = 'Guyana'
$ws.Cells.Item(154, 3).Value = 5
$ws.Cells.Item(154, 4).Value = 0
$ws.Cells.Item(154, 5).Value = 25
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 8).Value = 4

# Row 155
$ws.Cells.Item(155, 1).Value = 'Bahamas'
$ws.Cells.Item(155, 2).Value = 29
$ws.Cells.Item(155, 4).Value = 4
$ws.Cells.Item(155, 6).Value = 1
$ws.Cells.Item(155, 8).Value = 5

# Row 164
$ws.Cells.Item(164, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(164, 4).Value = 1
$ws.Cells.Item(164, 8).Value = 0

# Row 165
$ws.Cells.Item(165, 1).Value = 'Libia'
$ws.Cells.Item(165, 4).Value = 0
$ws.Cells.Item(165, 8).Value = 1

# Row 171
$ws.Cells.Item(171, 1).Value = 'Dominica'
$ws.Cells.Item(171, 3).Value = 0

# Row 172
$ws.Cells.Item(172, 1).Value = 'Fiyi'
$ws.Cells.Item(172, 3).Value = 2

# Row 181
$ws.Cells.Item(181, 1).Value = 'Seychelles'

# Row 182
$ws.Cells.Item(182, 1).Value = 'San Cristobal y Nieves'

# Row 185
$ws.Cells.Item(185, 1).Value = 'Suazilandia'

# Row 186
$ws.Cells.Item(186, 1).Value = 'Republica del Chad'

# Row 192
$ws.Cells.Item(192, 1).Value = 'Cabo Verde'
$ws.Cells.Item(192, 4).Value = 0
$ws.Cells.Item(192, 8).Value = 1

# Row 193
$ws.Cells.Item(193, 1).Value = 'Somalia'
$ws.Cells.Item(193, 4).Value = 1
$ws.Cells.Item(193, 8).Value = 0

# Row 194
$ws.Cells.Item(194, 1).Value = 'San Vicente y las Granadinas'

# Row 196
$ws.Cells.Item(196, 1).Value = 'San Bartolome'
$ws.Cells.Item(196, 4).Value = 1
$ws.Cells.Item(196, 8).Value = 0

# Row 197
$ws.Cells.Item(197, 1).Value = 'Botsuana'

# Row 198
$ws.Cells.Item(198, 1).Value = 'Nicaragua'
$ws.Cells.Item(198, 4).Value = 0
$ws.Cells.Item(198, 8).Value = 1

# Row 200
$ws.Cells.Item(200, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(200, 6).Value = 0
$ws.Cells.Item(200, 7).Value = 0

# Row 201
$ws.Cells.Item(201, 1).Value = 'Belice'
$ws.Cells.Item(201, 6).Value = 1
$ws.Cells.Item(201, 7).Value = 1

# Row 203
$ws.Cells.Item(203, 1).Value = 'Malaui'

# Row 204
$ws.Cells.Item(204, 1).Value = 'Sahara Occidental'

# Row 206
$ws.Cells.Item(206, 1).Value = 'Burundi'

# Row 207
$ws.Cells.Item(207, 1).Value = 'Anguila'

# Row 208
$ws.Cells.Item(208, 1).Value = 'Islas Virgenes Britanicas'

# Row 212
$ws.Cells.Item(212, 1).Value = 'Sudan del Sur'

# Row 213
$ws.Cells.Item(213, 1).Value = 'Papua Nueva Guinea'
